$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$new.Name = "Test1"
$wb.Styles.Item("Normal").Font.Size = 20
Write-Output $new.StandardHeight
